$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3529.75
$ws.Range("I86").Value = 4593.2856
$ws.Range("K86").Value = 4593.2856
$ws.Range("M86").Value = -3470.2856

$ws.Range("H89").Value = 3529.75
$ws.Range("I89").Value = 4593.2856
$ws.Range("K89").Value = 22966.428
$ws.Range("M89").Value = -17350.428

$ws.Range("H129").Value = 871.85364
$ws.Range("J129").Value = 893.8
$ws.Range("L129").Value = 2681.4
$ws.Range("N129").Value = -12681.4

$ws.Range("H132").Value = 10424682
$ws.Range("I132").Value = 11500659
$ws.Range("J132").Value = 23568.666
$ws.Range("K132").Value = 34501977
$ws.Range("L132").Value = 70705.99800000001
$ws.Range("M132").Value = -34499447
$ws.Range("N132").Value = -75765.99800000001

$ws.Range("H137").Value = 1815.8928
$ws.Range("I137").Value = 1245.55
$ws.Range("J137").Value = 3241.75
$ws.Range("K137").Value = 3736.65
$ws.Range("L137").Value = 9725.25
$ws.Range("M137").Value = -1186.65
$ws.Range("N137").Value = -14825.25

$ws.Range("H138").Value = 2800.561
$ws.Range("I138").Value = 2398.4546
$ws.Range("J138").Value = 2862.8591
$ws.Range("K138").Value = 7195.3638
$ws.Range("L138").Value = 8588.577300000001
$ws.Range("M138").Value = -2055.3638
$ws.Range("N138").Value = -18868.5773

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13203.725
$ws.Range("I32").Value = 10405.17
$ws.Range("J32").Value = 19182.455
$ws.Range("K32").Value = 10405.17
$ws.Range("L32").Value = 19182.455
$ws.Range("M32").Value = -10118.17
$ws.Range("N32").Value = -19756.455

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4692.7407
$ws.Range("I134").Value = 893.2632
$ws.Range("J134").Value = 13716.5
$ws.Range("K134").Value = 2679.7896
$ws.Range("L134").Value = 41149.5
$ws.Range("M134").Value = -144.7896000000001
$ws.Range("N134").Value = -46219.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -17246

$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -56232

$ws.Range("H105").Value = 710.6
$ws.Range("I105").Value = 706.2222
$ws.Range("J105").Value = 750
$ws.Range("K105").Value = 706.2222
$ws.Range("L105").Value = 750
$ws.Range("M105").Value = 1040.7778
$ws.Range("N105").Value = -4244

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4763990.5
$ws.Range("I34").Value = 387.16666
$ws.Range("J34").Value = 6669432
$ws.Range("K34").Value = 1161.49998
$ws.Range("L34").Value = 20008296
$ws.Range("M34").Value = -1077.49998
$ws.Range("N34").Value = -20008464

$ws.Range("H40").Value = 884.73334
$ws.Range("I40").Value = 67.42856999999999
$ws.Range("J40").Value = 1599.875
$ws.Range("K40").Value = 269.71428
$ws.Range("L40").Value = 6399.5
$ws.Range("M40").Value = -200.71428
$ws.Range("N40").Value = -6537.5

$ws.Range("H64").Value = 3453.2666
$ws.Range("I64").Value = 1333
$ws.Range("K64").Value = 3999
$ws.Range("M64").Value = -3729

$ws.Range("H67").Value = 3453.2666
$ws.Range("I67").Value = 1333
$ws.Range("K67").Value = 3999
$ws.Range("M67").Value = -3063

$ws.Range("H109").Value = 127252.125
$ws.Range("I109").Value = 251004.25
$ws.Range("J109").Value = 3500
$ws.Range("K109").Value = 753012.75
$ws.Range("L109").Value = 10500
$ws.Range("M109").Value = -751972.75
$ws.Range("N109").Value = -12580

$ws.Range("H117").Value = 787
$ws.Range("I117").Value = 641.875
$ws.Range("K117").Value = 1925.625
$ws.Range("M117").Value = 1516.375

$ws.Range("H137").Value = 28853002
$ws.Range("I137").Value = 62502410
$ws.Range("J137").Value = 10653.214
$ws.Range("K137").Value = 187507230
$ws.Range("L137").Value = 31959.642
$ws.Range("M137").Value = -187502130
$ws.Range("N137").Value = -42159.642

$ws.Range("H140").Value = 28085.975
$ws.Range("I140").Value = 47557.227
$ws.Range("J140").Value = 2887.8823
$ws.Range("K140").Value = 142671.681
$ws.Range("L140").Value = 8663.6469
$ws.Range("M140").Value = -137491.681
$ws.Range("N140").Value = -19023.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 15966.667
$ws.Range("J64").Value = 15966.667
$ws.Range("L64").Value = 15966.667
$ws.Range("N64").Value = -16462.667

$ws.Range("H67").Value = 15966.667
$ws.Range("J67").Value = 15966.667
$ws.Range("L67").Value = 15966.667
$ws.Range("N67").Value = -17682.667

$ws.Range("H86").Value = 32586.555
$ws.Range("J86").Value = 32586.555
$ws.Range("L86").Value = 32586.555
$ws.Range("N86").Value = -34958.555

$ws.Range("H89").Value = 32586.555
$ws.Range("J89").Value = 32586.555
$ws.Range("L89").Value = 97759.66500000001
$ws.Range("N89").Value = -109615.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1047.95
$ws.Range("I16").Value = 1047.95
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1047.95
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -877.95

$ws.Range("H68").Value = 1401
$ws.Range("I68").Value = 1401
$ws.Range("K68").Value = 1401
$ws.Range("M68").Value = -652

$ws.Range("H71").Value = 1401
$ws.Range("I71").Value = 1401
$ws.Range("K71").Value = 7005
$ws.Range("M71").Value = -3261

$ws.Range("H92").Value = 13500
$ws.Range("J92").Value = 13500
$ws.Range("L92").Value = 13500
$ws.Range("N92").Value = -18492

$ws.Range("H136").Value = 1689.381
$ws.Range("I136").Value = 1563.8823
$ws.Range("J136").Value = 2222.75
$ws.Range("K136").Value = 4691.6469
$ws.Range("L136").Value = 6668.25
$ws.Range("M136").Value = -2141.6469
$ws.Range("N136").Value = -11768.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 166668670

$ws.Range("H65").Value = 166668670

$ws.Range("H86").Value = 22325
$ws.Range("J86").Value = 22325
$ws.Range("L86").Value = 22325
$ws.Range("N86").Value = -24571

$ws.Range("H89").Value = 22325
$ws.Range("J89").Value = 22325
$ws.Range("L89").Value = 111625
$ws.Range("N89").Value = -122857

$ws.Range("H100").Value = 255.77777
$ws.Range("I100").Value = 255.77777
$ws.Range("K100").Value = 511.55554
$ws.Range("M100").Value = 29.44445999999999

$ws.Range("H113").Value = 249.36
$ws.Range("I113").Value = 149.94118
$ws.Range("J113").Value = 460.625
$ws.Range("K113").Value = 449.82354
$ws.Range("L113").Value = 1381.875
$ws.Range("M113").Value = 1720.17646
$ws.Range("N113").Value = -5721.875

$ws.Range("H136").Value = 1579.8864
$ws.Range("I136").Value = 675.0526
$ws.Range("J136").Value = 2267.56
$ws.Range("K136").Value = 2025.1578
$ws.Range("L136").Value = 6802.68
$ws.Range("M136").Value = 524.8422
$ws.Range("N136").Value = -11902.68
